# Uruguay Primera División - bases update (30-03-2024 19:32)
# Rows 118-120: opening-odds/Asian-handicap data rotated by one position
#   (the match previously stored in row120 now appears in row118, etc.)
# Rows 161-164: subsequent fixtures shifted up by two (new matches arrived,
#   two already-played/placeholder fixtures at the bottom - former rows
#   165/166 - are dropped) and rows 161/162 receive freshly-updated odds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 118 ---------------------------------------------------------
$ws.Range("B118").Value = 7013702
$ws.Range("F118").Value = 'Defensor Sporting'
$ws.Range("G118").Value = 'Danubio'
$ws.Range("K118").Value = 1.8
$ws.Range("L118").Value = 3.6
$ws.Range("M118").Value = 4.2
$ws.Range("N118").Value = 1.8
$ws.Range("O118").Value = 3.6
$ws.Range("P118").Value = 4.2
$ws.Range("Q118").Value = -0.75
$ws.Range("R118").Value = 2.05
$ws.Range("S118").Value = 1.8
$ws.Range("T118").Value = 2.25
$ws.Range("U118").Value = 1.85
$ws.Range("V118").Value = 2
$ws.Range("Y118").Value = 3.2
$ws.Range("AA118").Value = 0.8
$ws.Range("AB118").Value = -0.5
$ws.Range("AC118").Value = 0.5

# --- Row 119 ---------------------------------------------------------
$ws.Range("B119").Value = 7013885
$ws.Range("F119").Value = 'La Luz'
$ws.Range("G119").Value = 'Atletico Fenix Montevideo'
$ws.Range("I119").Value = 2
$ws.Range("K119").Value = 3
$ws.Range("L119").Value = 3
$ws.Range("M119").Value = 2.4
$ws.Range("N119").Value = 2.9
$ws.Range("O119").Value = 2.75
$ws.Range("P119").Value = 2.6
$ws.Range("Q119").Value = 0
$ws.Range("R119").Value = 2.025
$ws.Range("S119").Value = 1.825
$ws.Range("U119").Value = 2.025
$ws.Range("V119").Value = 1.825
$ws.Range("Y119").Value = 1.6
$ws.Range("AA119").Value = 0.825
$ws.Range("AB119").Value = 0
$ws.Range("AC119").Value = -0

# --- Row 120 ---------------------------------------------------------
$ws.Range("B120").Value = 7013886
$ws.Range("F120").Value = 'Racing Club de Montevideo'
$ws.Range("G120").Value = 'Cerro'
$ws.Range("I120").Value = 1
$ws.Range("K120").Value = 2.25
$ws.Range("L120").Value = 3.1
$ws.Range("M120").Value = 3.25
$ws.Range("N120").Value = 2.25
$ws.Range("O120").Value = 2.875
$ws.Range("P120").Value = 3.5
$ws.Range("Q120").Value = -0.25
$ws.Range("R120").Value = 1.95
$ws.Range("S120").Value = 1.9
$ws.Range("T120").Value = 2
$ws.Range("U120").Value = 1.925
$ws.Range("V120").Value = 1.925
$ws.Range("Y120").Value = 2.5
$ws.Range("AA120").Value = 0.8999999999999999
$ws.Range("AB120").Value = -1
$ws.Range("AC120").Value = 0.925

# --- Row 161 ---------------------------------------------------------
$ws.Range("B161").Value = 7994684
$ws.Range("E161").Value = 45381.79166666666
$ws.Range("F161").Value = 'Defensor Sporting'
$ws.Range("G161").Value = 'Danubio'
$ws.Range("K161").Value = 1.909
$ws.Range("L161").Value = 3.25
$ws.Range("M161").Value = 3.8
$ws.Range("N161").Value = 1.75
$ws.Range("O161").Value = 3.3
$ws.Range("P161").Value = 4.5
$ws.Range("Q161").Value = -0.75
$ws.Range("R161").Value = 2.025
$ws.Range("S161").Value = 1.825
$ws.Range("U161").Value = 1.975
$ws.Range("V161").Value = 1.875

# --- Row 162 ---------------------------------------------------------
$ws.Range("B162").Value = 7995146
$ws.Range("E162").Value = 45382.625
$ws.Range("F162").Value = 'Club Atletico Progreso'
$ws.Range("G162").Value = 'Deportivo Maldonado'
$ws.Range("K162").Value = 2.4
$ws.Range("L162").Value = 3.1
$ws.Range("M162").Value = 3
$ws.Range("N162").Value = 2.05
$ws.Range("O162").Value = 3.4
$ws.Range("P162").Value = 3.5
$ws.Range("Q162").Value = -0.25
$ws.Range("R162").Value = 1.8
$ws.Range("S162").Value = 2.05
$ws.Range("T162").Value = 2.5

# --- Row 163 ---------------------------------------------------------
$ws.Range("B163").Value = 7995141
$ws.Range("E163").Value = 45382.72916666666
$ws.Range("F163").Value = 'Miramar Misiones'
$ws.Range("G163").Value = 'Cerro'
$ws.Range("K163").Value = 2.6
$ws.Range("L163").Value = 3
$ws.Range("M163").Value = 2.75
$ws.Range("N163").Value = 2.5
$ws.Range("O163").Value = 3
$ws.Range("P163").Value = 2.9
$ws.Range("Q163").Value = 0
$ws.Range("R163").Value = 1.775
$ws.Range("S163").Value = 2.1
$ws.Range("U163").Value = 1.975
$ws.Range("V163").Value = 1.875

# --- Row 164 ---------------------------------------------------------
$ws.Range("B164").Value = 7994683
$ws.Range("E164").Value = 45382.83333333334
$ws.Range("F164").Value = 'Montevideo Wanderers'
$ws.Range("G164").Value = 'Boston River'
$ws.Range("K164").Value = 2.5
$ws.Range("M164").Value = 2.75
$ws.Range("N164").Value = 2.7
$ws.Range("O164").Value = 3
$ws.Range("P164").Value = 2.6
$ws.Range("Q164").Value = 0
$ws.Range("R164").Value = 1.975
$ws.Range("S164").Value = 1.875
$ws.Range("T164").Value = 2.25
$ws.Range("U164").Value = 2
$ws.Range("V164").Value = 1.85

# --- Drop the two trailing fixtures (former rows 165 & 166) ----------
$ws.Rows("165:166").Delete()
